$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20, column A previously pointed at shared string "Exp 23" (reused by mistake);
# correct it to reference "Exp 24"
$ws.Range("A20").Value = "Exp 24"

# Append new experiment row 21 with "Exp 25" parameters
$ws.Range("A21").Value = "Exp 25"
$ws.Range("B21").Value = 0.15
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "Local"
$ws.Range("E21").Value = -1
$ws.Range("F21").Value = "Exp 25.png"

# Match formatting (center alignment) used by the rest of the data rows (A:E)
$ws.Range("A21:E21").HorizontalAlignment = -4108

# Update the active selection to mirror the saved workbook state
$ws.Range("F18").Select()
